# Auto-generated update of cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{Row=2; D='43.535.75'; E='  -3.73%  '}
    @{Row=3; D='2.238.59'; E='  -4.90%  '}
    @{Row=4; E='  +0.08%  '}
    @{Row=5; D='318.38'; E='  +3.03%  '}
    @{Row=6; D='101.26'; E='  -7.41%  '}
    @{Row=7; D='0.584'; E='  -6.95%  '}
    @{Row=8; E='  +0.03%  '}
    @{Row=9; D='0.566'; E='  -7.94%  '}
    @{Row=10; D='37.16'; E='  -9.52%  '}
    @{Row=11; D='54.07'; E='  -2.29%  '}
    @{Row=12; D='0.0831'; E='  -9.07%  '}
    @{Row=13; D='7.77'; E='  -7.67%  '}
    @{Row=14; E='  -2.92%  '}
    @{Row=15; D='0.870'; E='  -11.39%  '}
    @{Row=16; D='2.577.95'; E='  -5.03%  '}
    @{Row=17; D='14.24'; E='  -7.10%  '}
    @{Row=18; D='2.236.48'; E='  -4.79%  '}
    @{Row=19; D='43.250.56'; E='  -4.36%  '}
    @{Row=20; D='14.25'; E='  +8.24%  '}
    @{Row=21; D='0.0₃0977'; E='  -8.21%  '}
    @{Row=22; D='6.59'; E='  -9.33%  '}
    @{Row=23; D='65.91'; E='  -10.06%  '}
    @{Row=24; D='3.23'; E='  -5.89%  '}
    @{Row=25; D='238.38'; E='  -7.84%  '}
    @{Row=26; D='2.17'; E='  -4.63%  '}
    @{Row=27; E='  -0.13%  '}
    @{Row=28; D='10.28'; E='  -6.85%  '}
    @{Row=29; E='  -6.47%  '}
    @{Row=30; D='6.52'; E='  -11.00%  '}
    @{Row=31; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0898'; E='  -6.68%  '}
    @{Row=32; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='20.73'; E='  -7.12%  '}
    @{Row=33; D='34.18'; E='  -9.36%  '}
    @{Row=34; D='157.74'; E='  -7.62%  '}
    @{Row=35; D='2.77'; E='  -5.93%  '}
    @{Row=36; D='3.22'; E='  +8.83%  '}
    @{Row=37; D='0.123'; E='  -6.02%  '}
    @{Row=38; D='4.57'; E='  -5.07%  '}
    @{Row=39; D='1.91'; E='  +11.21%  '}
    @{Row=40; D='0.104'; E='  -9.15%  '}
    @{Row=41; D='3.60'; E='  -8.18%  '}
    @{Row=42; D='0.0326'; E='  -8.52%  '}
    @{Row=43; D='1.00'; E='  +0.07%  '}
    @{Row=44; D='1.817.31'; E='  +11.86%  '}
    @{Row=45; B='Celestia'; C='https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; D='12.25'; E='  -3.91%  '}
    @{Row=46; B='BitcoinSV'; C='https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'; D='90.02'; E='  -9.02%  '}
    @{Row=47; B='ordi'; C='https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; D='80.07'; E='  -1.75%  '}
    @{Row=48; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.210'; E='  -9.42%  '}
    @{Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='16.63'; E='  +66.86%  '}
    @{Row=50; B='THORChain'; C='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; D='5.44'; E='  -1.36%  '}
    @{Row=51; B='MultiversX'; C='https://coinranking.com/coin/omwkOTglq+multiversx-egld'; D='61.41'; E='  -11.75%  '}
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) {
        $ws.Range("D$r").NumberFormat = '@'
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Range("E$r").NumberFormat = '@'
        $ws.Range("E$r").Value = $u.E
    }
}

Write-Host "Updated $($rowUpdates.Count) rows"
